# The AHB-Diff header row currently carries "<column>_old" / "<column>_new"
# suffixes (e.g. "Segmentname_old", "Segmentname_new"). Rename them to the
# concrete format-version suffixes "_FV2310" (old/before) and "_FV2404"
# (new/after). Column K stays "diff" and is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headersFV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}

$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# Turn the whole data range into a real Excel Table ("Table1") with an
# autofilter on the header row, mirroring xl/tables/table1.xml + the
# <tableParts> link added to the worksheet.
$rng = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split right below row 1, same as selecting
# A2 and choosing View > Freeze Panes > Freeze Top Row in the UI).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
